# "Update phần test case" — mark a batch of Test-Case (column K) statuses
# on the CODE-TC-SRS sheet as Done / In Progress (they were all "Not Start"
# before), matching the status-colour already used on the same row's
# SRS (I) / Code (J) status cells, then leave the workbook's view state
# (active sheet / selection) the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CODE-TC-SRS")

function Set-Status($row, $status, $styleSourceAddress) {
    $target = $ws.Range("K$row")
    $target.Value = $status
    $ws.Range($styleSourceAddress).Copy()
    $target.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = 0
}

Set-Status 22 "Done"         "J22"
Set-Status 23 "In Progress"  "J23"
Set-Status 24 "Done"         "J24"
Set-Status 25 "Done"         "J25"
Set-Status 26 "Done"         "J26"
Set-Status 27 "Done"         "J27"
Set-Status 28 "In Progress"  "I28"
Set-Status 42 "Done"         "J42"
Set-Status 43 "Done"         "J43"
Set-Status 44 "Done"         "J44"
Set-Status 45 "In Progress"  "I45"
Set-Status 46 "In Progress"  "K45"

# Restore the view state: TỔNG HỢP selection at G5, then DOC (not selected
# on return), finally CODE-TC-SRS active with K46 selected/scrolled into
# view - matching where the author ended up after making the edits above.
$wsSummary = $wb.Worksheets.Item("TỔNG HỢP")
$wsSummary.Activate()
$wsSummary.Range("G5").Select()

$wsDoc = $wb.Worksheets.Item("DOC")
$wsDoc.Activate()

$ws.Activate()
$ws.Range("K46").Select()
